$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, shifting the existing rows 87-93 down to 88-94.
$ws.Rows("87:87").Insert()

# Populate the newly inserted row with the weekly price entry (2023-10-13).
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 45212
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112012
$ws.Range("G87").Value = "Espinaca"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 35
$ws.Range("K87").Value = 12000
$ws.Range("L87").Value = 12000
$ws.Range("M87").Value = 12000
$ws.Range("N87").Value = "`$/cuna 10 kilos"
$ws.Range("O87").Value = "Región Metropolitana"
$ws.Range("P87").Value = 1200
$ws.Range("Q87").Value = 10
$ws.Range("R87").Value = "Hortaliza"
